$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 15927400
$ws.Range("E8").Value = 13641600
$ws.Range("F8").Value = 13561500
$ws.Range("G8").Value = 13279600
$ws.Range("H8").Value = 12862100
$ws.Range("I8").Value = 12867300
$ws.Range("J8").Value = 12991300

# Row 17
$ws.Range("D17").Value = 6629300
$ws.Range("E17").Value = 5955900
$ws.Range("F17").Value = 4642200
$ws.Range("G17").Value = 3018300
$ws.Range("H17").Value = 2489000
$ws.Range("I17").Value = 4997300
$ws.Range("J17").Value = 3552000

# Row 18
$ws.Range("D18").Value = 9298100
$ws.Range("E18").Value = 7685700
$ws.Range("F18").Value = 8919300
$ws.Range("G18").Value = 10261300
$ws.Range("H18").Value = 10373100
$ws.Range("I18").Value = 7870000
$ws.Range("J18").Value = 9439300

# Row 20
$ws.Range("D20").Value = -1710900
$ws.Range("E20").Value = -3343100
$ws.Range("F20").Value = 1898000
$ws.Range("G20").Value = -1305900
$ws.Range("H20").Value = -3806900
$ws.Range("I20").Value = 132000
$ws.Range("J20").Value = -3447300

# Row 21
$ws.Range("D21").Value = 9208600
$ws.Range("E21").Value = 5892000
$ws.Range("F21").Value = 12337100
$ws.Range("G21").Value = 10409100
$ws.Range("H21").Value = 8025500
$ws.Range("I21").Value = 9424400
$ws.Range("J21").Value = 7491800

# Row 23
$ws.Range("D23").Value = 7587300
$ws.Range("E23").Value = 4342600
$ws.Range("F23").Value = 10817300
$ws.Range("G23").Value = 8955300
$ws.Range("H23").Value = 6566100
$ws.Range("I23").Value = 8002000
$ws.Range("J23").Value = 5992000

# Row 24
$ws.Range("D24").Value = 2147900
$ws.Range("E24").Value = 824800
$ws.Range("F24").Value = 3132700
$ws.Range("G24").Value = 2757100
$ws.Range("H24").Value = 2044000
$ws.Range("I24").Value = 36400
$ws.Range("J24").Value = 125500

# Row 26
$ws.Range("D26").Value = 5439300
$ws.Range("E26").Value = 3517700
$ws.Range("F26").Value = 7684600
$ws.Range("G26").Value = 6198200
$ws.Range("H26").Value = 4522100
$ws.Range("I26").Value = 7965700
$ws.Range("J26").Value = 5866600

# Row 27
$ws.Range("D27").Value = 5221600
$ws.Range("E27").Value = 3276500
$ws.Range("F27").Value = 7666500
$ws.Range("G27").Value = 5487500
$ws.Range("H27").Value = 4445300
$ws.Range("I27").Value = 7839400
$ws.Range("J27").Value = 5855400

# Row 32
$ws.Range("D32").Value = 1710900
$ws.Range("E32").Value = 3343100
$ws.Range("F32").Value = -1898000
$ws.Range("G32").Value = 1305900
$ws.Range("H32").Value = 3806900
$ws.Range("I32").Value = -132000
$ws.Range("J32").Value = 3447300

# Row 33
$ws.Range("D33").Value = 5221600
$ws.Range("E33").Value = 3276500
$ws.Range("F33").Value = 7666500
$ws.Range("G33").Value = 5487500
$ws.Range("H33").Value = 4445300
$ws.Range("I33").Value = 7839400
$ws.Range("J33").Value = 5855400

# Row 35
$ws.Range("D35").Value = 5221600
$ws.Range("E35").Value = 3276500
$ws.Range("F35").Value = 7666500
$ws.Range("G35").Value = 5487500
$ws.Range("H35").Value = 4445300
$ws.Range("I35").Value = 7839400
$ws.Range("J35").Value = 5855400

# Row 41
$ws.Range("D41").Value = 435465000
$ws.Range("E41").Value = 430186000
$ws.Range("F41").Value = 331316000
$ws.Range("G41").Value = 263029000
$ws.Range("H41").Value = 15339800
$ws.Range("I41").Value = 11466700
$ws.Range("J41").Value = 10998300

# Row 42
$ws.Range("D42").Value = 359968000
$ws.Range("E42").Value = 368742000
$ws.Range("F42").Value = 403950000
$ws.Range("G42").Value = 291768000
$ws.Range("H42").Value = 549634000
$ws.Range("I42").Value = 550699000
$ws.Range("J42").Value = 464937000

# Row 47
$ws.Range("G47").Value = 2639400
$ws.Range("H47").Value = 2410200
$ws.Range("I47").Value = 2725300
$ws.Range("J47").Value = 3815500

# Row 48
$ws.Range("D48").Value = 5020200
$ws.Range("E48").Value = 5776700
$ws.Range("F48").Value = 16615400
$ws.Range("G48").Value = 9745600
$ws.Range("H48").Value = 12263600
$ws.Range("I48").Value = 9871600
$ws.Range("J48").Value = 9988000

# Row 49
$ws.Range("D49").Value = 15734000
$ws.Range("E49").Value = 14387900
$ws.Range("F49").Value = 612400
$ws.Range("G49").Value = 5944300
$ws.Range("H49").Value = 637300
$ws.Range("I49").Value = 638400
$ws.Range("J49").Value = 686500

# Row 52
$ws.Range("D52").Value = 516100
$ws.Range("E52").Value = 574200
$ws.Range("F52").Value = 518400
$ws.Range("G52").Value = 7054100
$ws.Range("H52").Value = 3662300
$ws.Range("I52").Value = 5801700
$ws.Range("J52").Value = 7564000

# Row 54
$ws.Range("D54").Value = 1846470000
$ws.Range("E54").Value = 1812130000
$ws.Range("F54").Value = 1752040000
$ws.Range("G54").Value = 1714750000
$ws.Range("H54").Value = 1588320000
$ws.Range("I54").Value = 1615870000
$ws.Range("J54").Value = 1503910000

# Row 57
$ws.Range("D57").Value = 16873800
$ws.Range("E57").Value = 17012800
$ws.Range("F57").Value = 30059600

# Row 59
$ws.Range("D59").Value = 2688400
$ws.Range("E59").Value = 2557500
$ws.Range("F59").Value = 2514500
$ws.Range("G59").Value = 541200
$ws.Range("H59").Value = 1963800
$ws.Range("I59").Value = 1773800
$ws.Range("J59").Value = 1902700

# Row 61
$ws.Range("D61").Value = 117115000
$ws.Range("E61").Value = 131346000
$ws.Range("F61").Value = 133480000
$ws.Range("G61").Value = 119415000
$ws.Range("H61").Value = 89079600
$ws.Range("I61").Value = 79572100
$ws.Range("J61").Value = 76494800

# Row 62
$ws.Range("D62").Value = 2768100
$ws.Range("E62").Value = 1269900
$ws.Range("F62").Value = 1824800
$ws.Range("G62").Value = 6508400
$ws.Range("H62").Value = 300000
$ws.Range("I62").Value = 129300
$ws.Range("J62").Value = 143300

# Row 66
$ws.Range("D66").Value = 1766300000
$ws.Range("E66").Value = 1737440000
$ws.Range("F66").Value = 1679590000
$ws.Range("G66").Value = 1640940000
$ws.Range("H66").Value = 1530660000
$ws.Range("I66").Value = 1564090000
$ws.Range("J66").Value = 1463490000

# Row 70
$ws.Range("F70").Value = 894300
$ws.Range("G70").Value = 1926600
$ws.Range("H70").Value = 2826400
$ws.Range("I70").Value = 3411300
$ws.Range("J70").Value = 3709700

# Row 72
$ws.Range("D72").Value = 11807500
$ws.Range("E72").Value = 8306800
$ws.Range("F72").Value = 6750900
$ws.Range("G72").Value = 25069600
$ws.Range("H72").Value = -4858800
$ws.Range("I72").Value = -7985800
$ws.Range("J72").Value = -14519200

# Row 76
$ws.Range("D76").Value = 80170500
$ws.Range("E76").Value = 74682700
$ws.Range("F76").Value = 71557200
$ws.Range("G76").Value = 71884500
$ws.Range("H76").Value = 54835000
$ws.Range("I76").Value = 48370900
$ws.Range("J76").Value = 36706000

# Row 81
$ws.Range("D81").Value = 5221600
$ws.Range("E81").Value = 3276500
$ws.Range("F81").Value = 7666500
$ws.Range("G81").Value = 5487500
$ws.Range("H81").Value = 4445300
$ws.Range("I81").Value = 7839400
$ws.Range("J81").Value = 5855400

# Row 83
$ws.Range("D83").Value = 1619600
$ws.Range("E83").Value = 1547700
$ws.Range("F83").Value = 1518100
$ws.Range("G83").Value = 1452200
$ws.Range("H83").Value = 1457800
$ws.Range("I83").Value = 1420800
$ws.Range("J83").Value = 1498100

# Row 89
$ws.Range("D89").Value = -1425300
$ws.Range("E89").Value = 10968300
$ws.Range("F89").Value = 2185100
$ws.Range("G89").Value = 60160800
$ws.Range("H89").Value = 53799800
$ws.Range("I89").Value = -11733800
$ws.Range("J89").Value = -16176500

# Row 91
$ws.Range("D91").Value = -2641500
$ws.Range("E91").Value = -3895500
$ws.Range("F91").Value = -3512600
$ws.Range("G91").Value = -1955300
$ws.Range("H91").Value = -4131100
$ws.Range("I91").Value = -1458700
$ws.Range("J91").Value = -1577600

# Row 94
$ws.Range("D94").Value = -39662800
$ws.Range("E94").Value = -115726000
$ws.Range("F94").Value = -54041800
$ws.Range("G94").Value = 23677800
$ws.Range("H94").Value = 3772400
$ws.Range("I94").Value = -57149400
$ws.Range("J94").Value = -16659800

# Row 96
$ws.Range("D96").Value = -1721100
$ws.Range("E96").Value = -1717900
$ws.Range("F96").Value = -1765400
$ws.Range("G96").Value = -1592700
$ws.Range("H96").Value = -1375600
$ws.Range("I96").Value = -1378700
$ws.Range("J96").Value = -1951800

# Row 100
$ws.Range("D100").Value = 42059700
$ws.Range("E100").Value = 107302600
$ws.Range("F100").Value = 50250800
$ws.Range("G100").Value = -8166700
$ws.Range("H100").Value = -53986900
$ws.Range("I100").Value = 69064300
$ws.Range("J100").Value = 26811200

# Row 101
$ws.Range("D101").Value = -125200
$ws.Range("E101").Value = -107500
$ws.Range("F101").Value = -253700
$ws.Range("G101").Value = 339600
$ws.Range("H101").Value = 287800
$ws.Range("I101").Value = 287300
$ws.Range("J101").Value = -12900

# Row 102
$ws.Range("D102").Value = 846500
$ws.Range("E102").Value = 2437100
$ws.Range("F102").Value = -1859600
$ws.Range("G102").Value = 76011500
$ws.Range("H102").Value = 3873100
$ws.Range("I102").Value = 468400
$ws.Range("J102").Value = -6037900

